$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.946.26'
$ws.Range("E2").Value = '  +1.34%  '
$ws.Range("D3").Value = '3.183.33'
$ws.Range("E3").Value = '  +0.73%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'595.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.59%  '
$ws.Range("D6").Value = "'153.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.11%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.181.03'
$ws.Range("E8").Value = '  +0.70%  '
$ws.Range("D9").Value = "'0.539"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.19%  '
$ws.Range("D10").Value = "'0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.06%  '
$ws.Range("D11").Value = "'6.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("D12").Value = "'0.515"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.47%  '
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("D14").Value = "'38.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.79%  '
$ws.Range("D15").Value = '3.704.32'
$ws.Range("E15").Value = '  +0.70%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '66.026.81'
$ws.Range("E16").Value = '  +1.41%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").Value = "'7.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.75%  '
$ws.Range("D18").Value = '3.181.55'
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("E20").Value = '  +0.61%  '
$ws.Range("D21").Value = "'15.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.44%  '
$ws.Range("E22").Value = '  +2.64%  '
$ws.Range("E23").Value = '  +4.07%  '
$ws.Range("D24").Value = "'15.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.81%  '
$ws.Range("D25").Value = "'84.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.65%  '
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").Value = "'9.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.54%  '
$ws.Range("D28").Value = "'2.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.63%  '
$ws.Range("E29").Value = '  +5.34%  '
$ws.Range("D30").Value = "'6.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +12.87%  '
$ws.Range("D31").Value = "'2.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.66%  '
$ws.Range("D32").Value = "'28.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.84%  '
$ws.Range("E33").Value = '  +2.27%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("D35").Value = "'6.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.72%  '
$ws.Range("D36").Value = "'54.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("D37").Value = "'0.0901"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").Value = "'481.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.00%  '
$ws.Range("D39").Value = "'0.0419"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("E40").Value = '  +1.61%  '
$ws.Range("D41").Value = "'0.301"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.35%  '
$ws.Range("E42").Value = '  +3.59%  '
$ws.Range("E43").Value = '  -5.19%  '
$ws.Range("D44").Value = '0.0₃0653'
$ws.Range("E44").Value = '  +10.86%  '
$ws.Range("D45").Value = '2.895.52'
$ws.Range("E45").Value = '  -4.92%  '
$ws.Range("D46").Value = "'2.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").Value = "'28.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.24%  '
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("E49").Value = '  +1.54%  '
$ws.Range("E50").Value = '  +2.74%  '
$ws.Range("D51").Value = "'2.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.49%  '
